$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price entry was added for "Papa" at Feria Lagunitas de Puerto Montt.
# It belongs at row 630 (sheet rows are sorted with newest entries per-variety group
# interleaved), pushing every existing row from 630 downward by one.
$ws.Rows.Item(630).Insert()

# Populate the newly inserted row with the new observation's data.
$ws.Cells.Item(630, 1).Value = 4
$ws.Cells.Item(630, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(630, 3).Value = "Los Lagos"
$ws.Cells.Item(630, 4).Value = 45127
$ws.Cells.Item(630, 5).Value = 10
$ws.Cells.Item(630, 6).Value = 100114001
$ws.Cells.Item(630, 7).Value = "Papa"
$ws.Cells.Item(630, 8).Value = "Patagonia"
$ws.Cells.Item(630, 9).Value = "1a (guarda)"
$ws.Cells.Item(630, 10).Value = 250
$ws.Cells.Item(630, 11).Value = 18000
$ws.Cells.Item(630, 12).Value = 18000
$ws.Cells.Item(630, 13).Value = 18000
$ws.Cells.Item(630, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(630, 15).Value = "Provincia de Llanquihue"
$ws.Cells.Item(630, 16).Value = 720
$ws.Cells.Item(630, 17).Value = 25
$ws.Cells.Item(630, 18).Value = "Hortaliza"
